# Add reports for carry increment adder
# Also tidy up a couple of nearby precision-only values and fill in
# previously-missing data points (row 17: H17 / O17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 (Verilog ('+') version of adders block) -----------------------
$ws.Range("G7").Value = 333.8
$ws.Range("U7").Value = 34.62

# --- Row 12 (Carry Increment adder block, columns J:O) -------------------
# plus one cell in the adjoining "Ripple Carry Adder" block (H12) and one
# in the "Carry Select Adder" block (V12).
$ws.Range("H12").Value = 181.8

$ws.Range("J12").Value = 204.0
$ws.Range("K12").Value = 2656.4
$ws.Range("L12").Value = 16843.6
$ws.Range("M12").Value = 14187.0
$ws.Range("N12").Value = 41.5
$ws.Range("O12").Value = 2656.4

$ws.Range("V12").Value = 67.6

# --- Row 17 (Carry Save / Carry Skip Adder block) -------------------------
# H17 previously held the shared text "20ns"; replace with the real number.
$ws.Range("H17").Value = 156.7
# O17 was previously blank; fill in the reading.
$ws.Range("O17").Value = 155.7
